$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.507.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "'1.666.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'233.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D7").Value = "'0.4612"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.94%  "
$ws.Range("D8").Value = "'0.2571"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'0.06117"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "'1.664.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").Value = "'0.06954"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "'14.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'4.329"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "'74.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").Value = "'0.5649"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.50%  "
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "'25.498.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.21%  "
$ws.Range("D19").Value = "'0.000006665"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").Value = "'1.880.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "'8.688"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "'5.195"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "'136.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").Value = "'14.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "'1.366"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").Value = "'103.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("D30").Value = "'3.940"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").Value = "'0.07720"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").Value = "'3.591"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'0.04267"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("E34").Value = "  +1.78%  "
$ws.Range("D35").Value = "'0.9399"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.5939"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.29%  "
$ws.Range("D37").Value = "'0.9145"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.11%  "
$ws.Range("D38").Value = "'2.503"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "'1.001"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").Value = "'101.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("D41").Value = "'0.01459"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.50%  "
$ws.Range("D42").Value = "'1.805"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("D43").Value = "'0.3686"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "'4.904"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.40%  "
$ws.Range("D45").Value = "'0.05254"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").Value = "'0.1100"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").Value = "'6.102"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("D48").Value = "'29.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("B49").Value = "TrueUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D49").Value = "'1.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "'1.001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.321"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.39%  "
